$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "runs"/"balls"/"fours"/"sixes" figures (columns C-F) are stored as
# text in this sheet, so force a text format before writing the new
# values back in order to keep them as text (matching the rest of the
# sheet) instead of letting them be auto-detected as numbers.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2
Set-TextValue "C2" "17"
Set-TextValue "D2" "12"
Set-TextValue "F2" "2"

# Row 3
Set-TextValue "C3" "4"
Set-TextValue "D3" "7"

# Row 4
Set-TextValue "C4" "0"
Set-TextValue "D4" "1"
Set-TextValue "F4" "0"

# Row 5
Set-TextValue "C5" "51"
Set-TextValue "D5" "26"
Set-TextValue "E5" "6"
Set-TextValue "F5" "1"

# Row 6
Set-TextValue "C6" "16"
Set-TextValue "D6" "18"
Set-TextValue "E6" "1"

# Row 8
Set-TextValue "C8" "7"
Set-TextValue "D8" "14"
Set-TextValue "E8" "0"
Set-TextValue "F8" "0"

# Row 9
Set-TextValue "C9" "3"
Set-TextValue "D9" "5"

# Row 10
Set-TextValue "C10" "12"
Set-TextValue "D10" "13"
Set-TextValue "E10" "1"

# Row 11
Set-TextValue "C11" "8"
Set-TextValue "D11" "7"
Set-TextValue "E11" "0"
